$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.598495
$ws.Range("H2").Value = 61.795485
$ws.Range("I2").Value = 0.187290770808477
$ws.Range("J2").Value = 0.187290770808477
$ws.Range("M2").Value = 7.655977
$ws.Range("N2").Value = 22.967931
$ws.Range("O2").Value = 0.2994795900616967
$ws.Range("P2").Value = 0.2994795900616967
$ws.Range("Q2").Value = 157.701603954615
$ws.Range("R2").Value = 1419.314435591535
$ws.Range("S2").Value = 0.05608976326406188
$ws.Range("T2").Value = 0.05608976326406188
$ws.Range("G3").Value = 20.598495
$ws.Range("H3").Value = 61.795485
$ws.Range("I3").Value = 0.187290770808477
$ws.Range("J3").Value = 0.187290770808477
$ws.Range("O3").Value = 0.3140620915319453
$ws.Range("P3").Value = 0.3140620915319453
$ws.Range("Q3").Value = 165.38053750416
$ws.Range("R3").Value = 1488.42483753744
$ws.Range("S3").Value = 0.0588209312047405
$ws.Range("T3").Value = 0.0588209312047405
$ws.Range("G4").Value = 20.598495
$ws.Range("H4").Value = 61.795485
$ws.Range("I4").Value = 0.187290770808477
$ws.Range("J4").Value = 0.187290770808477
$ws.Range("M4").Value = 9.879524666666667
$ws.Range("N4").Value = 29.638574
$ws.Range("O4").Value = 0.386458318406358
$ws.Range("P4").Value = 0.386458318406358
$ws.Range("Q4").Value = 203.50333944871
$ws.Range("R4").Value = 1831.53005503839
$ws.Range("S4").Value = 0.07238007633967464
$ws.Range("T4").Value = 0.07238007633967464
$ws.Range("G5").Value = 60.20577233333334
$ws.Range("I5").Value = 0.5474179306512287
$ws.Range("J5").Value = 0.5474179306512288
$ws.Range("M5").Value = 7.655977
$ws.Range("N5").Value = 22.967931
$ws.Range("O5").Value = 0.2994795900616967
$ws.Range("P5").Value = 0.2994795900616967
$ws.Range("Q5").Value = 460.9340082512364
$ws.Range("R5").Value = 4148.406074261127
$ws.Range("S5").Value = 0.1639404974638523
$ws.Range("T5").Value = 0.1639404974638523
$ws.Range("G6").Value = 60.20577233333334
$ws.Range("I6").Value = 0.5474179306512287
$ws.Range("J6").Value = 0.5474179306512288
$ws.Range("O6").Value = 0.3140620915319453
$ws.Range("P6").Value = 0.3140620915319453
$ws.Range("Q6").Value = 483.378178325152
$ws.Range("S6").Value = 0.1719232202424143
$ws.Range("T6").Value = 0.1719232202424143
$ws.Range("G7").Value = 60.20577233333334
$ws.Range("I7").Value = 0.5474179306512287
$ws.Range("J7").Value = 0.5474179306512288
$ws.Range("M7").Value = 9.879524666666667
$ws.Range("N7").Value = 29.638574
$ws.Range("O7").Value = 0.386458318406358
$ws.Range("P7").Value = 0.386458318406358
$ws.Range("Q7").Value = 594.8044128428843
$ws.Range("R7").Value = 5353.239715585958
$ws.Range("S7").Value = 0.2115542129449622
$ws.Range("T7").Value = 0.2115542129449622
$ws.Range("G8").Value = 29.17709966666666
$ws.Range("H8").Value = 87.53129899999999
$ws.Range("I8").Value = 0.2652912985402942
$ws.Range("J8").Value = 0.2652912985402942
$ws.Range("M8").Value = 7.655977
$ws.Range("N8").Value = 22.967931
$ws.Range("O8").Value = 0.2994795900616967
$ws.Range("P8").Value = 0.2994795900616967
$ws.Range("Q8").Value = 223.3792039747077
$ws.Range("R8").Value = 2010.412835772369
$ws.Range("S8").Value = 0.0794493293337825
$ws.Range("T8").Value = 0.0794493293337825
$ws.Range("G9").Value = 29.17709966666666
$ws.Range("H9").Value = 87.53129899999999
$ws.Range("I9").Value = 0.2652912985402942
$ws.Range("J9").Value = 0.2652912985402942
$ws.Range("O9").Value = 0.3140620915319453
$ws.Range("P9").Value = 0.3140620915319453
$ws.Range("Q9").Value = 234.256164136544
$ws.Range("R9").Value = 2108.305477228896
$ws.Range("S9").Value = 0.08331794008479051
$ws.Range("T9").Value = 0.08331794008479051
$ws.Range("G10").Value = 29.17709966666666
$ws.Range("H10").Value = 87.53129899999999
$ws.Range("I10").Value = 0.2652912985402942
$ws.Range("J10").Value = 0.2652912985402942
$ws.Range("M10").Value = 9.879524666666667
$ws.Range("N10").Value = 29.638574
$ws.Range("O10").Value = 0.386458318406358
$ws.Range("P10").Value = 0.386458318406358
$ws.Range("Q10").Value = 288.2558758586251
$ws.Range("R10").Value = 2594.302882727626
$ws.Range("S10").Value = 0.1025240291217212
$ws.Range("T10").Value = 0.1025240291217212
